$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object 'object[,]' 24,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.049185476250953
$arrBF[0,2] = 1.056631418948794
$arrBF[0,3] = 1.056424711825644
$arrBF[0,4] = 1.06713958029171
$arrBF[1,0] = 1.02
$arrBF[1,1] = 1.050042291345221
$arrBF[1,2] = 1.05730876261192
$arrBF[1,3] = 1.057169986226655
$arrBF[1,4] = 1.067924808290429
$arrBF[2,0] = 1.02
$arrBF[2,1] = 1.050597475527165
$arrBF[2,2] = 1.057747720460268
$arrBF[2,3] = 1.057653250772433
$arrBF[2,4] = 1.06843393389966
$arrBF[3,0] = 1.02
$arrBF[3,1] = 1.050831056983225
$arrBF[3,2] = 1.057932417176784
$arrBF[3,3] = 1.057856657988262
$arrBF[3,4] = 1.068648215034502
$arrBF[4,0] = 1.02
$arrBF[4,1] = 1.050870286961657
$arrBF[4,2] = 1.057963437828637
$arrBF[4,3] = 1.057890825137249
$arrBF[4,4] = 1.068684208072889
$arrBF[5,0] = 1.02
$arrBF[5,1] = 1.050600595940924
$arrBF[5,2] = 1.057750187765039
$arrBF[5,3] = 1.057655967757966
$arrBF[5,4] = 1.068436796175563
$arrBF[6,0] = 1.02
$arrBF[6,1] = 1.049474880967544
$arrBF[6,2] = 1.056860190320937
$arrBF[6,3] = 1.056676368003242
$arrBF[6,4] = 1.067404737180444
$arrBF[7,0] = 1.02
$arrBF[7,1] = 1.047497181886414
$arrBF[7,2] = 1.055297123603673
$arrBF[7,3] = 1.05495810116187
$arrBF[7,4] = 1.065594098735263
$arrBF[8,0] = 1.02
$arrBF[8,1] = 1.046182822056968
$arrBF[8,2] = 1.054258701051067
$arrBF[8,3] = 1.053818024382452
$arrBF[8,4] = 1.064392491979943
$arrBF[9,0] = 1.02
$arrBF[9,1] = 1.04561468533382
$arrBF[9,2] = 1.053809935113368
$arrBF[9,3] = 1.053325671697133
$arrBF[9,4] = 1.063873509350491
$arrBF[10,0] = 1.02
$arrBF[10,1] = 1.045403804501021
$arrBF[10,2] = 1.05364337725666
$arrBF[10,3] = 1.053142988479437
$arrBF[10,4] = 1.063680936618225
$arrBF[11,0] = 1.02
$arrBF[11,1] = 1.045449032319646
$arrBF[11,2] = 1.053679098397337
$arrBF[11,3] = 1.053182165636035
$arrBF[11,4] = 1.063722235015882
$arrBF[12,0] = 1.02
$arrBF[12,1] = 1.045597250779321
$arrBF[12,2] = 1.053796164648624
$arrBF[12,3] = 1.05331056698483
$arrBF[12,4] = 1.063857587125858
$arrBF[13,0] = 1.02
$arrBF[13,1] = 1.04568859305144
$arrBF[13,2] = 1.053868310815574
$arrBF[13,3] = 1.053389705654476
$arrBF[13,4] = 1.063941008660824
$arrBF[14,0] = 1.02
$arrBF[14,1] = 1.046220548421869
$arrBF[14,2] = 1.054288502826717
$arrBF[14,3] = 1.053850727936395
$arrBF[14,4] = 1.064426963155744
$arrBF[15,0] = 1.02
$arrBF[15,1] = 1.046554496186545
$arrBF[15,2] = 1.054552314550764
$arrBF[15,3] = 1.054140266633813
$arrBF[15,4] = 1.064732144698517
$arrBF[16,0] = 1.02
$arrBF[16,1] = 1.046749377647782
$arrBF[16,2] = 1.054706276041267
$arrBF[16,3] = 1.054309275736672
$arrBF[16,4] = 1.06491027925115
$arrBF[17,0] = 1.02
$arrBF[17,1] = 1.046815843315027
$arrBF[17,2] = 1.054758787218335
$arrBF[17,3] = 1.054366924798533
$arrBF[17,4] = 1.064971040063678
$arrBF[18,0] = 1.02
$arrBF[18,1] = 1.04651865687367
$arrBF[18,2] = 1.054524001287733
$arrBF[18,3] = 1.054109188816592
$arrBF[18,4] = 1.064699388422855
$arrBF[19,0] = 1.02
$arrBF[19,1] = 1.045553599970231
$arrBF[19,2] = 1.053761687846785
$arrBF[19,3] = 1.053272750492801
$arrBF[19,4] = 1.063817723751092
$arrBF[20,0] = 1.02
$arrBF[20,1] = 1.044947701671429
$arrBF[20,2] = 1.053283166160862
$arrBF[20,3] = 1.052747997904439
$arrBF[20,4] = 1.063264547365341
$arrBF[21,0] = 1.02
$arrBF[21,1] = 1.045268816653364
$arrBF[21,2] = 1.053536765461752
$arrBF[21,3] = 1.053026069575468
$arrBF[21,4] = 1.063557685880894
$arrBF[22,0] = 1.02
$arrBF[22,1] = 1.046534850821979
$arrBF[22,2] = 1.054536794571493
$arrBF[22,3] = 1.05412323115516
$arrBF[22,4] = 1.064714189180715
$arrBF[23,0] = 1.02
$arrBF[23,1] = 1.048007748290789
$arrBF[23,2] = 1.055700583522109
$arrBF[23,3] = 1.055401364820919
$arrBF[23,4] = 1.066061234043675
$ws.Range("B2:F25").Value = $arrBF

$arrIM = New-Object 'object[,]' 24,5
$arrIM[0,0] = 1.047076762066251
$arrIM[0,1] = 1.054225262611738
$arrIM[0,2] = 1.059368599651256
$arrIM[0,3] = 1.05916245964086
$arrIM[0,4] = 1.069848244470518
$arrIM[1,0] = 1.047287393483905
$arrIM[1,1] = 1.054731710113886
$arrIM[1,2] = 1.059859967496469
$arrIM[1,3] = 1.059721544458103
$arrIM[1,4] = 1.070449276200888
$arrIM[2,0] = 1.047422918653766
$arrIM[2,1] = 1.055059501569578
$arrIM[2,2] = 1.060177903540011
$arrIM[2,3] = 1.060083662331434
$arrIM[2,4] = 1.070838551887947
$arrIM[3,0] = 1.047479708944336
$arrIM[3,1] = 1.055197324348303
$arrIM[3,2] = 1.060311559632247
$arrIM[3,3] = 1.060235979608429
$arrIM[3,4] = 1.071002289862425
$arrIM[4,0] = 1.047489233443219
$arrIM[4,1] = 1.0552204664889
$arrIM[4,2] = 1.060334000790259
$arrIM[4,3] = 1.060261559168707
$arrIM[4,4] = 1.071029787212796
$arrIM[5,0] = 1.047423678214283
$arrIM[5,1] = 1.055061343090344
$arrIM[5,2] = 1.060179689477588
$arrIM[5,3] = 1.060085697278387
$arrIM[5,4] = 1.070840739425015
$arrIM[6,0] = 1.047148104384164
$arrIM[6,1] = 1.05439640028493
$arrIM[6,2] = 1.059534661322374
$arrIM[6,3] = 1.05935133106428
$arrIM[6,4] = 1.07005128879567
$arrIM[7,0] = 1.046656664256108
$arrIM[7,1] = 1.053225405002875
$arrIM[7,2] = 1.058398010516439
$arrIM[7,3] = 1.058060050600894
$arrIM[7,4] = 1.068663065641188
$arrIM[8,0] = 1.046325154177575
$arrIM[8,1] = 1.05244530556842
$arrIM[8,2] = 1.057640307237372
$arrIM[8,3] = 1.057201146295875
$arrIM[8,4] = 1.06773962304555
$arrIM[9,0] = 1.046180695238632
$arrIM[9,1] = 1.052107664754305
$arrIM[9,2] = 1.057312246121028
$arrIM[9,3] = 1.056829712473977
$arrIM[9,4] = 1.067340266390971
$arrIM[10,0] = 1.046126900373312
$arrIM[10,1] = 1.051982273441457
$arrIM[10,2] = 1.057190395453165
$arrIM[10,3] = 1.056691818667107
$arrIM[10,4] = 1.067192004334723
$arrIM[11,0] = 1.046138445713821
$arrIM[11,1] = 1.052009169217339
$arrIM[11,2] = 1.05721653254808
$arrIM[11,3] = 1.056721394017415
$arrIM[11,4] = 1.067223803547751
$arrIM[12,0] = 1.046176251316014
$arrIM[12,1] = 1.052097299381531
$arrIM[12,2] = 1.057302173779509
$arrIM[12,3] = 1.05681831262572
$arrIM[12,4] = 1.067328009427157
$arrIM[13,0] = 1.046199526546148
$arrIM[13,1] = 1.052151602441999
$arrIM[13,2] = 1.057354940984742
$arrIM[13,3] = 1.056878037131383
$arrIM[13,4] = 1.067392224333154
$arrIM[14,0] = 1.046334722273577
$arrIM[14,1] = 1.052467716899967
$arrIM[14,2] = 1.057662080313931
$arrIM[14,3] = 1.057225807324195
$arrIM[14,4] = 1.067766137708896
$arrIM[15,0] = 1.046419283116915
$arrIM[15,1] = 1.052666047578747
$arrIM[15,2] = 1.057854749575138
$arrIM[15,3] = 1.057444083272555
$arrIM[15,4] = 1.068000818745465
$arrIM[16,0] = 1.046468517894188
$arrIM[16,1] = 1.052781744636795
$arrIM[16,2] = 1.057967133037246
$arrIM[16,3] = 1.057571445901033
$arrIM[16,4] = 1.068137752363402
$arrIM[17,0] = 1.046485290708862
$arrIM[17,1] = 1.052821196698332
$arrIM[17,2] = 1.05800545331543
$arrIM[17,3] = 1.057614881015274
$arrIM[17,4] = 1.068184451303439
$arrIM[18,0] = 1.046410219649073
$arrIM[18,1] = 1.052644767102589
$arrIM[18,2] = 1.05783407769298
$arrIM[18,3] = 1.05742065956409
$arrIM[18,4] = 1.067975634699981
$arrIM[19,0] = 1.046165122267122
$arrIM[19,1] = 1.052071346584475
$arrIM[19,2] = 1.05727695439885
$arrIM[19,3] = 1.056789770476812
$arrIM[19,4] = 1.06729732125996
$arrIM[20,0] = 1.046010231124711
$arrIM[20,1] = 1.051710951100655
$arrIM[20,2] = 1.056926703349914
$arrIM[20,3] = 1.056393530040969
$arrIM[20,4] = 1.066871283760792
$arrIM[21,0] = 1.046092416351781
$arrIM[21,1] = 1.05190199014327
$arrIM[21,2] = 1.057112374315466
$arrIM[21,3] = 1.05660354383857
$arrIM[21,4] = 1.067097091558757
$arrIM[22,0] = 1.046414315312468
$arrIM[22,1] = 1.052654382790277
$arrIM[22,2] = 1.057843418418426
$arrIM[22,3] = 1.057431243587869
$arrIM[22,4] = 1.06798701413757
$arrIM[23,0] = 1.046784401128701
$arrIM[23,1] = 1.053528041842383
$arrIM[23,2] = 1.058691856827515
$arrIM[23,3] = 1.058393540500242
$arrIM[23,4] = 1.069021602199975
$ws.Range("I2:M25").Value = $arrIM

Write-Output "Updated vm_pu values for case with 380 kV"
